# Only strip punctuation if it occurs at the end of a word
# Fix test data: the "expected longest word" for the apostrophe/hyphen
# examples should retain the apostrophe/hyphen since it's in the middle
# of the word, not at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: "longest word is apostrophized"
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = "Shouldn't"

# Row 7: "longest word is hyphenated"
$ws.Range("C7").Value = 19
$ws.Range("D7").Value = "properly-hyphenated"

# Update view state to match final selection
$ws.Range("D7").Select()
